# Apply edits: add a "methodology" column (classification/regression) after
# "Algorithm" and before "Initial Balance", and expand the data from 6 rows
# (TSLA/AAPL x Ensemble/LSTM/VWAP) to 10 rows (TSLA/AAPL x
# Ensemble/LSTM/VWAP x classification/regression, with VWAP only appearing
# once per ticker, without a methodology value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Initial Balance") so the existing
# Initial Balance / Final Balance / Cumulative Returns columns shift right
# by one (C->D, D->E, E->F), preserving their header style/formatting.
$ws.Range("C1").EntireColumn.Insert()

# New header cell for the inserted column, matching the style used by the
# other header cells (bold, bordered, centered - same as existing row 1).
$ws.Range("C1").Value = "methodology"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Clear any old data below the header so we can rewrite the full table
# cleanly (old sheet had 7 rows total, new has 11).
$ws.Range("A2:F100").ClearContents()

# Full data set (Stock Ticker, Algorithm, methodology, Initial Balance,
# Final Balance, Cumulative Returns (%))
$data = @(
    @("TSLA", "Ensemble", "classification", 1000, 985.8895705521473,  -1.411042944785268),
    @("TSLA", "LSTM",     "classification", 1000, 896.7208884810494,  -10.32791115189505),
    @("TSLA", "Ensemble", "regression",     1000, 894.8501091417539,  -10.51498908582461),
    @("TSLA", "LSTM",     "regression",     1000, 1001.870779339296,  0.1870779339295581),
    @("TSLA", "VWAP",     "",               1000, 1000.570825680124,  0.05708256801236758),
    @("AAPL", "Ensemble", "classification", 1000, 1016.698409886167,  1.669840988616733),
    @("AAPL", "LSTM",     "classification", 1000, 1016.698409886167,  1.669840988616733),
    @("AAPL", "Ensemble", "regression",     1000, 991.3163165158546,  -0.8683683484145388),
    @("AAPL", "LSTM",     "regression",     1000, 1000,               0),
    @("AAPL", "VWAP",     "",               1000, 1124.453764535177,  12.4453764535177)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne "") {
        $ws.Cells.Item($r, 3).Value = $row[2]
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}
